# "semana 30 de 2025" - add week-30 column (AG) to the weekly IRA hospital
# report, and correct week-28's (AE) existing tally on row 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell AG1 = "30" -------------------------------------------
# This must be stored as TEXT (like the other week-number headers D1:AF1,
# which are all inline/shared strings), reusing the existing bold+centered
# header style (s="1") without minting a brand-new style record (and
# without Excel's usual "number-stored-as-text" quote-prefix marking, since
# the source headers were never authored that way). Routing the literal
# number through TEXT() and collapsing the formula to its cached value via
# Copy/PasteSpecial(Values) keeps both the string type and the original
# style intact.
$ws.Range("AG1").Formula = "=TEXT(30,""0"")"
$ws.Range("AG1").Copy() | Out-Null
$ws.Range("AG1").PasteSpecial(-4163) | Out-Null

# --- Correction on existing cell AE28: 163 -> 74 ---------------------------
$ws.Range("AE28").Value = 74

# --- New week-30 numbers (column AG) for every row that already carries an
# explicit (non-omitted) value through column AF. A handful of rows (3, 11,
# 13, 15, 18, 19, 20, 21, 24, 26, 27, 33, 51) are sparse in the source data
# (their trailing zero weeks are simply absent) and stay untouched.
$ws.Range("AG2").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AG5").Value = 0
$ws.Range("AG6").Value = 23
$ws.Range("AG7").Value = 2
$ws.Range("AG8").Value = 28
$ws.Range("AG9").Value = 0
$ws.Range("AG10").Value = 0
$ws.Range("AG12").Value = 0
$ws.Range("AG14").Value = 0
$ws.Range("AG16").Value = 0
$ws.Range("AG17").Value = 0
$ws.Range("AG22").Value = 0
$ws.Range("AG23").Value = 0
$ws.Range("AG25").Value = 2
$ws.Range("AG28").Value = 44
$ws.Range("AG29").Value = 4
$ws.Range("AG30").Value = 18
$ws.Range("AG31").Value = 0
$ws.Range("AG32").Value = 0
$ws.Range("AG34").Value = 5
$ws.Range("AG35").Value = 20
$ws.Range("AG36").Value = 0
$ws.Range("AG37").Value = 0
$ws.Range("AG38").Value = 0
$ws.Range("AG39").Value = 0
$ws.Range("AG40").Value = 0
$ws.Range("AG41").Value = 0
$ws.Range("AG42").Value = 0
$ws.Range("AG43").Value = 0
$ws.Range("AG44").Value = 0
$ws.Range("AG45").Value = 0
$ws.Range("AG46").Value = 0
$ws.Range("AG47").Value = 0
$ws.Range("AG48").Value = 0
$ws.Range("AG49").Value = 0
$ws.Range("AG50").Value = 0
$ws.Range("AG52").Value = 0
$ws.Range("AG53").Value = 0
$ws.Range("AG54").Value = 0
$ws.Range("AG55").Value = 0
$ws.Range("AG56").Value = 0
$ws.Range("AG57").Value = 0
